# Bulk data upload (Employee information) bux fixes
#
# Adds a new "UserId" column (G) to the employee template sheet, mirroring
# the formatting of the existing header cells, and moves the active
# selection to G3 (just below the new header, as in the authored edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell in column G, continuing the existing header row.
$ws.Range("G1").Value = "UserId"

# Match the bold formatting used by the other header cells (A1:F1).
$ws.Range("G1").Font.Bold = $true

# Update the active selection to G3, as captured in the saved workbook.
$ws.Range("G3").Select()
